$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Print"
$ws.Range("A6").Value = 1
$ws.Range("A6").NumberFormat = "0.00"

$ws.Range("C7").Select()

$ws.PageSetup.Orientation = $excel.xlPortrait
